$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 5, 6, 7 are cyclically rotated:
#   new row5 <- old row7
#   new row6 <- old row5
#   new row7 <- old row6
# Only columns A, B, D, E, F, G, H, I, Q, R, AC carry row-specific data that
# differs between the three rows; capture those before overwriting.

$cols = @("A","B","D","E","F","G","H","I","Q","R","AC")

# Column I is stored as text (e.g. "1") rather than as a number in this
# workbook; force the text number format on it so re-assigning the
# captured value doesn't get silently re-coerced back into a numeric cell.
$ws.Range("I5:I7").NumberFormat = "@"

$old5 = @{}
$old6 = @{}
$old7 = @{}
foreach ($col in $cols) {
    $old5[$col] = $ws.Range($col + "5").Value2
    $old6[$col] = $ws.Range($col + "6").Value2
    $old7[$col] = $ws.Range($col + "7").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "5").Value2 = $old7[$col]
    $ws.Range($col + "6").Value2 = $old5[$col]
    $ws.Range($col + "7").Value2 = $old6[$col]
}
